$p = $ppt.ActivePresentation

# --- 1. Slide master title placeholder: collapse the per-letter runs into
#        the normal placeholder prompt text "Click to edit Master title style".
$master = $p.SlideMaster
$title = $master.Shapes.Item("PlaceHolder 1")
$title.TextFrame.TextRange.Text = "Click to edit Master title style"

# --- 2. Slide 2: widen/reposition the two red "rightArrow" connector shapes.

$s2 = $p.Slides.Item(2)

# "CustomShape 15" - keep its right-hand tip fixed, extend the tail to the
# left (left edge moves from 482.4pt to 446.4pt, width grows 64.8pt -> 100.8pt).
# (values nudged by <1/1000 pt so the float32 EMU round-trip lands on the
# exact target EMU instead of 1 EMU short)
$arrow15 = $s2.Shapes.Item("CustomShape 15")
$arrow15.Left = 446.40000999999995
$arrow15.Top = 223.20000499999998
$arrow15.Width = 100.8
$arrow15.Height = 25.483464566929133

# "CustomShape 56" - left edge stays put, arrow just gets longer
# (width grows 64.8pt -> 108pt).
$arrow56 = $s2.Shapes.Item("CustomShape 56")
$arrow56.Left = 158.400002
$arrow56.Top = 223.20000499999998
$arrow56.Width = 108.0
$arrow56.Height = 25.483464566929133
